# Update each stock's worksheet with the latest remn_amt figures:
#  - fill in the previously-placeholder (0) values for 2025-10-29 / 2025-10-30
#  - append a new row for 2025-10-31 (date serial 45961) with remn_amt = 0
#
# Values per sheet (in tab order), taken from the source upload:
#   LG생활건강   : B100 -> 471234, B101 -> 464924
#   아모레퍼시픽 : B100 -> 357883, B101 -> 364003
#   한국콜마     : B100 -> 181933, B101 -> 178109
#   코스맥스     : B100 -> 257372, B101 -> 255552
#   에이피알     : B100 -> 530238, B101 -> 563157
#   달바글로벌   : B100 -> 58212,  B101 -> 57920

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = 1; B100 = 471234; B101 = 464924 },
    @{ Sheet = 2; B100 = 357883; B101 = 364003 },
    @{ Sheet = 3; B100 = 181933; B101 = 178109 },
    @{ Sheet = 4; B100 = 257372; B101 = 255552 },
    @{ Sheet = 5; B100 = 530238; B101 = 563157 },
    @{ Sheet = 6; B100 = 58212;  B101 = 57920 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)

    # Fill in the two trailing values that were still placeholders.
    $ws.Cells.Item(100, 2).Value = $u.B100
    $ws.Cells.Item(101, 2).Value = $u.B101

    # Append the next date row (2025-10-31), matching the date format used
    # by the existing date column (A), with remn_amt left at 0.
    $dateCell = $ws.Cells.Item(102, 1)
    $dateCell.Value = 45961
    $dateCell.NumberFormat = $ws.Cells.Item(101, 1).NumberFormat
    $ws.Cells.Item(102, 2).Value = 0
}
